$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 1.53
$ws.Range("Q2").Value = 1.58

# Row 6
$ws.Range("Q6").Value = 1.68
$ws.Range("R6").Value = 2.05

# Row 8
$ws.Range("G8").Value = 1.72
$ws.Range("Q8").Value = 1.77

# Row 9
$ws.Range("G9").Value = 3.3
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 1.92
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 2.3
$ws.Range("L9").Value = 2.6
$ws.Range("Q9").Value = 1.63
$ws.Range("R9").Value = 2.15
$ws.Range("X9").Value = 19
$ws.Range("AA9").Value = 23
$ws.Range("AG9").Value = 151
$ws.Range("AJ9").Value = 9
$ws.Range("AK9").Value = 19
$ws.Range("AN9").Value = 5.5
$ws.Range("AO9").Value = 17
$ws.Range("AX9").Value = 11
$ws.Range("AY9").Value = 19
$ws.Range("AZ9").Value = 34
$ws.Range("BA9").Value = 51

# Row 10
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 3.25
$ws.Range("J10").Value = 2.63
$ws.Range("K10").Value = 2.38
$ws.Range("L10").Value = 3.6
$ws.Range("W10").Value = 10
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 9
$ws.Range("Z10").Value = 19
$ws.Range("AA10").Value = 15
$ws.Range("AE10").Value = 12
$ws.Range("AG10").Value = 126
$ws.Range("AL10").Value = 23
$ws.Range("AN10").Value = 4.33
$ws.Range("AO10").Value = 11
$ws.Range("AP10").Value = 17
$ws.Range("AQ10").Value = 34
$ws.Range("AR10").Value = 41
$ws.Range("AX10").Value = 17
$ws.Range("BC10").Value = 351

# Row 12
$ws.Range("K12").Value = 2.4
$ws.Range("L12").Value = 2.63
$ws.Range("N12").Value = 17
$ws.Range("O12").Value = 1.14
$ws.Range("P12").Value = 5.5
$ws.Range("Q12").Value = 1.53
$ws.Range("R12").Value = 2.4
$ws.Range("U12").Value = 1.44
$ws.Range("V12").Value = 2.63
$ws.Range("AB12").Value = 23
$ws.Range("AC12").Value = 17
$ws.Range("AF12").Value = 29
$ws.Range("AH12").Value = 12
$ws.Range("AJ12").Value = 9.5
$ws.Range("AP12").Value = 19
$ws.Range("AQ12").Value = 41
$ws.Range("BB12").Value = 81
$ws.Range("BC12").Value = 301

# Row 13
$ws.Range("G13").Value = 1.38
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 7.5
$ws.Range("L13").Value = 9
$ws.Range("Z13").Value = 8.5
$ws.Range("AH13").Value = 13
$ws.Range("AJ13").Value = 26
$ws.Range("AM13").Value = 81
$ws.Range("AN13").Value = 3.1
$ws.Range("AQ13").Value = 23

# Row 14
$ws.Range("K14").Value = 1.91

# Row 15
$ws.Range("A15").Value = "65Ly84z2"
$ws.Range("C15").Value = "17:15"
$ws.Range("D15").Value = "PERU - LIGA 1"
$ws.Range("E15").Value = "AD Tarma"
$ws.Range("F15").Value = "Union Comercio"
$ws.Range("G15").Value = 1.14
$ws.Range("H15").Value = 7.5
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 1.44
$ws.Range("K15").Value = 3.2
$ws.Range("L15").Value = 13
$ws.Range("M15").Value = 1.02
$ws.Range("N15").Value = 21
$ws.Range("O15").Value = 1.11
$ws.Range("P15").Value = 6.5
$ws.Range("Q15").Value = 1.4
$ws.Range("R15").Value = 2.88
$ws.Range("S15").Value = 1.2
$ws.Range("T15").Value = 4.33
$ws.Range("U15").Value = 2.1
$ws.Range("V15").Value = 1.67
$ws.Range("W15").Value = 10
$ws.Range("X15").Value = 7
$ws.Range("Y15").Value = 11
$ws.Range("Z15").Value = 7
$ws.Range("AA15").Value = 11
$ws.Range("AC15").Value = 21
$ws.Range("AD15").Value = 15
$ws.Range("AE15").Value = 29
$ws.Range("AF15").Value = 67
$ws.Range("AH15").Value = 34
$ws.Range("AI15").Value = 67
$ws.Range("AJ15").Value = 34
$ws.Range("AK15").Value = 201
$ws.Range("AL15").Value = 81
$ws.Range("AM15").Value = 67
$ws.Range("AN15").Value = 3.25
$ws.Range("AO15").Value = 4.75
$ws.Range("AP15").Value = 17
$ws.Range("AQ15").Value = 9.5
$ws.Range("AR15").Value = 29
$ws.Range("AS15").Value = 101
$ws.Range("AT15").Value = 4.33
$ws.Range("AU15").Value = 11
$ws.Range("AV15").Value = 67
$ws.Range("AW15").Value = 15
$ws.Range("AX15").Value = 51
$ws.Range("AY15").Value = 51
$ws.Range("AZ15").Value = 351
$ws.Range("BA15").Value = 301
$ws.Range("BB15").Value = 351
$ws.Range("BC15").Value = 151
$ws.Range("BD15").Value = 151

# Row 16
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 2.75
$ws.Range("K16").Value = 2.05
$ws.Range("Q16").Value = 2.25
$ws.Range("R16").Value = 1.58
$ws.Range("X16").Value = 9
$ws.Range("Y16").Value = 9
$ws.Range("Z16").Value = 17
$ws.Range("AH16").Value = 10
$ws.Range("AI16").Value = 19
$ws.Range("AJ16").Value = 15
$ws.Range("AP16").Value = 23
$ws.Range("AZ16").Value = 81

# Row 17
$ws.Range("R17").Value = 1.54

# Row 18
$ws.Range("R18").Value = 1.44

# Row 19
$ws.Range("G19").Value = 1.85
$ws.Range("Q19").Value = 1.69
$ws.Range("R19").Value = 2.07
